$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sample data row (row 2) first - removes the old demo strings.
$ws.Rows.Item(2).Delete()

# --- Re-shape the header row: insert new columns so surviving headers land
# in their final positions (mirrors an author inserting columns next to the
# related existing header). ---

# Gender / Nationality before the old "Special Service Request Code" column (D).
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# Baggage Weight / Origin before the old "Destination" column (now H).
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(8).Insert()

# Is VIP? before the old "Phone Number" column (now L).
$ws.Columns.Item(12).Insert()

# --- Fill in the header text, in the order the shared-string table records. ---

# Columns that keep their existing text (already correct, reused as-is).
$ws.Range("B1").Value = "Passenger's Name"
$ws.Range("C1").Value = "Passenger's Seat Number"
$ws.Range("F1").Value = "Special Service Request Code"
$ws.Range("K1").Value = "PNR"
$ws.Range("J1").Value = "Passsenger's Destination"
$ws.Range("G1").Value = "Passenger's Baggage Count(Please enter digits)"
$ws.Range("N1").Value = "PAX Type "

# Brand-new trailing columns (O:S).
$ws.Range("O1").Value = "Identification Document Type"
$ws.Range("P1").Value = "Identification Document Number"
$ws.Range("Q1").Value = "Inbound Flight Number"
$ws.Range("R1").Value = "Outbound Flight Number"

# Newly inserted columns D:E.
$ws.Range("D1").Value = "Passenger's Gender"
$ws.Range("E1").Value = "Passenger's Nationality"

# Column A: "Flight Number" becomes "Affected Flight Number".
$ws.Range("A1").Value = "Affected Flight Number"

# Newly inserted columns H:I.
$ws.Range("H1").Value = "Baggage Weight(In KG)"
$ws.Range("I1").Value = "Passsenger's Origin"

# Final trailing column S.
$ws.Range("S1").Value = "Employee Id(If Passenger is an employee)"

# Columns O:S sit past the sheet's old used range, so they need the
# wrap-text header style applied explicitly (columns A:N inherited it already).
$ws.Range("O1:S1").WrapText = $true

# Column M: " Phone Number" (leading space) becomes "Phone Number".
$ws.Range("M1").Value = "Phone Number"

# Newly inserted column L.
$ws.Range("L1").Value = 'Is VIP? (Please Enter "Yes" or "No")'

# --- Styling / layout ---
$ws.Rows.Item(1).RowHeight = 60

$ws.Columns.Item(4).ColumnWidth = 11.833333333333334
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(8).ColumnWidth = 13.333333333333334
$ws.Columns.Item(9).ColumnWidth = 13.333333333333334
$ws.Columns.Item(10).ColumnWidth = 12.166666666666666
$ws.Columns.Item(12).ColumnWidth = 11.666666666666666
$ws.Columns.Item(15).ColumnWidth = 14.333333333333334
$ws.Columns.Item(16).ColumnWidth = 17.666666666666668
$ws.Columns.Item(17).ColumnWidth = 18.0
$ws.Columns.Item(18).ColumnWidth = 14.166666666666666
$ws.Columns.Item(19).ColumnWidth = 13.666666666666666

$ws.Range("A7").Select()
$ws.Rows.Item(7).Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
